$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = 42611.884664351855
$ws.Range("B7").Value = -8
$ws.Range("C7").Value = 47
$ws.Range("D7").Value = 50
$ws.Range("E7").Value = 40
$ws.Range("F7").Value = 60
$ws.Range("G7").Value = 10999
$ws.Range("H7").Value = 9577
$ws.Range("I7").Value = 1590
$ws.Range("J7").Value = 117
$ws.Range("K7").Value = 124
$ws.Range("L7").Value = 4
$ws.Range("M7").Value = 6
$ws.Range("N7").Value = "Noun"
